$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'74.634.48"
$ws.Range("E2").Value = "  +9.06%  "

$ws.Range("D3").Value = "'2.592.06"
$ws.Range("E3").Value = "  +6.62%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'186.53"
$ws.Range("E5").Value = "  +16.39%  "

$ws.Range("D6").Value = "'578.61"
$ws.Range("E6").Value = "  +3.74%  "

$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").Value = "'0.537"
$ws.Range("E8").Value = "  +5.61%  "

$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "'0.207"
$ws.Range("E9").Value = "  +28.00%  "

$ws.Range("D10").Value = "'2.591.08"
$ws.Range("E10").Value = "  +6.66%  "

$ws.Range("E11").Value = "  -0.21%  "

$ws.Range("E12").Value = "  +8.88%  "

$ws.Range("D13").Value = "'4.80"
$ws.Range("E13").Value = "  +3.43%  "

$ws.Range("E14").Value = "  +11.46%  "

$ws.Range("D15").Value = "'74.507.15"
$ws.Range("E15").Value = "  +9.11%  "

$ws.Range("D16").Value = "'3.057.91"
$ws.Range("E16").Value = "  +6.42%  "

$ws.Range("E17").Value = "  +13.85%  "

$ws.Range("D18").Value = "'2.608.44"
$ws.Range("E18").Value = "  +7.53%  "

$ws.Range("D19").Value = "'8.48"
$ws.Range("E19").Value = "  +23.02%  "

$ws.Range("D20").Value = "'11.68"
$ws.Range("E20").Value = "  +11.67%  "

$ws.Range("D21").Value = "'377.70"
$ws.Range("E21").Value = "  +12.77%  "

$ws.Range("D22").Value = "'2.31"
$ws.Range("E22").Value = "  +21.06%  "

$ws.Range("E23").Value = "  +6.49%  "

$ws.Range("E24").Value = "  +0.15%  "

$ws.Range("D25").Value = "'69.95"
$ws.Range("E25").Value = "  +5.01%  "

$ws.Range("D26").Value = "'4.18"
$ws.Range("E26").Value = "  +13.82%  "

$ws.Range("D27").Value = "'9.19"
$ws.Range("E27").Value = "  +11.82%  "

$ws.Range("D28").Value = "'2.726.92"
$ws.Range("E28").Value = "  +6.89%  "

$ws.Range("E29").Value = "  -0.03%  "

$ws.Range("D30").Value = "'0.0₃0940"
$ws.Range("E30").Value = "  +15.27%  "

$ws.Range("E31").Value = "  +11.23%  "

$ws.Range("D32").Value = "'501.20"
$ws.Range("E32").Value = "  +17.79%  "

$ws.Range("D33").Value = "'1.35"
$ws.Range("E33").Value = "  +17.46%  "

$ws.Range("E34").Value = "  +6.08%  "

$ws.Range("E35").Value = "  +0.06%  "

$ws.Range("E36").Value = "  +14.29%  "

$ws.Range("D37").Value = "'159.99"
$ws.Range("E37").Value = "  +0.74%  "

$ws.Range("D38").Value = "'19.19"
$ws.Range("E38").Value = "  +7.28%  "

$ws.Range("D39").Value = "'19.40"
$ws.Range("E39").Value = "  +1.90%  "

$ws.Range("E40").Value = "  +0.04%  "

$ws.Range("D41").Value = "'4.98"
$ws.Range("E41").Value = "  +15.04%  "

$ws.Range("E42").Value = "  +12.54%  "

$ws.Range("D43").Value = "'0.319"
$ws.Range("E43").Value = "  +7.50%  "

$ws.Range("E44").Value = "  +20.26%  "

$ws.Range("D45").Value = "'39.14"
$ws.Range("E45").Value = "  +4.95%  "

$ws.Range("E46").Value = "  +7.64%  "

$ws.Range("D47").Value = "'148.14"
$ws.Range("E47").Value = "  +11.32%  "

$ws.Range("D48").Value = "'0.0814"
$ws.Range("E48").Value = "  +14.21%  "

$ws.Range("E49").Value = "  +8.58%  "

$ws.Range("D50").Value = "'0.519"
$ws.Range("E50").Value = "  +8.00%  "

$ws.Range("D51").Value = "'0.581"
$ws.Range("E51").Value = "  +4.60%  "
